# JViewer.pptx maintenance update.
#
# The underlying commit is a housekeeping pass (Loaddb.ijs retrieval
# back-off, Run.ijs working-directory change, version bump) that happened
# to be saved from PowerPoint on 2024-03-07. Opening/saving the deck
# re-stamped every "today's date" footer field (the ppPlaceholderDate
# placeholder on the slide master and on each of its layouts) from the
# previous save date (12/19/23) to the new one (3/7/24). That is the only
# durable content change in the package - a slide was briefly added,
# tweaked, and then deleted again during the same session, so it leaves no
# trace in the final slide deck.
#
# Walk the slide master and every custom (slide) layout, find the date
# placeholder on each, and stamp the new date over its displayed text.

$p = $ppt.ActivePresentation

$newDate = "3/7/24"
$ppPlaceholderDate = 16

function Update-DatePlaceholder {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)

        $phType = $null
        try { $phType = $shape.PlaceholderFormat.Type } catch { $phType = $null }

        if ($phType -eq $ppPlaceholderDate -and $shape.HasTextFrame) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master footer date.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's footer date.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-DatePlaceholder $layouts.Item($li).Shapes
}
